$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Invoice date moved from 2021-06-21 to 2021-06-28
$ws.Range("F2").Value = "2021-06-28 "

# Invoice # changed from 3 to 6 (numeric-looking text -> keep as text with leading apostrophe)
$ws.Range("F3").Value = "'6"

# Bill-to customer block (name / address / phone) now populated so that
# loyalty points earned get credited to the right customer account
$ws.Range("A9").Value = "Sant Anurag Deo"
$ws.Range("A11").Value = "102, Whitestone veroso, Banglore 49"
$ws.Range("A13").Value = "'9900019362"

# Salesperson table: customer name + contact number
$ws.Range("B16").Value = "Sant Anurag Deo"
$ws.Range("D16").Value = "'9900019362"

# Line item updated: different book, qty and unit price
$ws.Range("B19").Value = "Test Book106"
$ws.Range("C19").Value = "'1"
$ws.Range("D19").Value = "'100"
